$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the generic category hints with specific per-question hints.
# Values are assigned in this particular order so that the workbook's
# shared-string table gets new entries appended in the same order as the
# target file (Calendar, Time, History, Numbers, Great Outdoors, Money,
# Measurements, Poetry, Space, English, Weather).
$ws.Range("B4").Value = "Calendar"          # 365 Days IN A Year

$ws.Range("B8").Value = "Time"              # 52 Weeks IN A Year
$ws.Range("B18").Value = "Time"             # 86400 Seconds IN A Day
$ws.Range("B26").Value = "Time"             # 10 Years IN A Decade

$ws.Range("B27").Value = "History"          # 45 Presidents IN America

$ws.Range("B28").Value = "Numbers"          # 12 Zeros IN A Billion

$ws.Range("B29").Value = "Great Outdoors"   # 40 National Parks IN Canada

$ws.Range("B22").Value = "Money"            # 40 Quarters IN A Roll

$ws.Range("B6").Value = "Measurements"      # 13 Items IN A Baker's Dozen
$ws.Range("B20").Value = "Measurements"     # 16 Tablespoons IN A Cup

$ws.Range("B19").Value = "Poetry"           # 5 Stanzas IN A Limerick

$ws.Range("B7").Value = "Space"             # 8 Planets IN THE Milky Way

$ws.Range("B2").Value = "English"           # 5 Vowels IN THE Alphabet
$ws.Range("B5").Value = "English"           # 26 Letters IN THE Alphabet

$ws.Range("B3").Value = "Weather"           # 7 Colors IN THE Rainbow

# Remaining rows keep an existing category (Geography / Entertainment /
# Sports), so no change is required for them.

# Remove the rows for questions that were dropped entirely. Deleted from
# bottom to top so earlier row numbers remain valid.
$ws.Rows(25).Delete()  # "118 Elements ON THE Periodic Table"
$ws.Rows(17).Delete()  # "7 Books IN THE Chronicles Of Narnia"
$ws.Rows(14).Delete()  # "4 Countries IN THE United Kingdom"

# Restore the selection to match the target workbook.
$ws.Range("A12").Select()
